$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Rebuild-Paragraph($matchText, $newInnerXml) {
    $r = $d.Content
    $find = $r.Find
    $find.ClearFormatting()
    $ok = $find.Execute($matchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Host "NOT FOUND: $matchText"
        return
    }
    $r.Delete()
    $frag = '<w:p ' + $wNs + '>' + $newInnerXml + '</w:p>'
    $r.InsertXML($frag)
}

# ---------------------------------------------------------------------------
# Hunk 1: reposition the gramStart proofErr so it wraps "disc[" instead of
# sitting between "disc" and "[" - for the 5 occurrences in the "disc[]"
# explanation section.
# ---------------------------------------------------------------------------

Rebuild-Paragraph "Ans:  disc[] = Discance of each node as we traverse DFS in the graph." (
    '<w:r><w:t>Ans:</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t>[</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve">] = </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Disc</w:t></w:r>' +
    '<w:r><w:t>ance</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> of each node as we traverse DFS in the graph.</w:t></w:r>'
)

Rebuild-Paragraph "   Increment disc[] value for each node as we detect the next node." (
    '<w:r><w:t xml:space="preserve">   Increment </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t>[</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>] value for each node as we detect the next node.</w:t></w:r>'
)

Rebuild-Paragraph "   node 1 -> disc[0]++ = 1" (
    '<w:r><w:t xml:space="preserve">   node 1 -&gt; </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t>[</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>0]++ = 1</w:t></w:r>'
)

Rebuild-Paragraph "   node 2 -> disc[1]++ = 2" (
    '<w:r><w:t xml:space="preserve">   node 2 -&gt; </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t>[</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>1]++ = 2</w:t></w:r>'
)

Rebuild-Paragraph "   node 3 -> disc[2]++ = 3" (
    '<w:r><w:t xml:space="preserve">   node 3 -&gt; </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t>[</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t>2]++ = 3</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# Hunk 2 (part 1): same gramStart reposition for the 6th occurrence, in the
# "Min[u] = min(min[u], disc[v])..." paragraph.
# ---------------------------------------------------------------------------

Rebuild-Paragraph "    Min[u] = min(min[u], disc[v]). Where disc[v] is the disc[ ] value of backedge and min[u] is the min[ ] value of current node whose min is getting calculated." (
    '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Min[u] = min(min[u], </w:t></w:r>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t>[v]).</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> Where </w:t></w:r>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">[v] is the </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>disc</w:t></w:r>' +
    '<w:r><w:t>[</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> ] value of </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>backedge</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> and min[u] is the min[ ] value of current node whose min is getting calculated.</w:t></w:r>'
)

# ---------------------------------------------------------------------------
# Hunk 2 (part 2): drop the old "_GoBack" bookmark that sits after
# "If min[u] >=", and append the new highlighted sentence (with the bookmark
# moved alongside it) to the end of the paragraph that now ends with
# "...min is getting calculated."
# ---------------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$r = $d.Content
$find = $r.Find
$find.ClearFormatting()
$found = $find.Execute("min is getting calculated.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $newTailXml = '<w:p ' + $wNs + '>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Make sure, v is not a parent of u but ancestor.</w:t></w:r>' +
        '</w:p>'
    $r.InsertXML($newTailXml)
} else {
    Write-Host "NOT FOUND: min is getting calculated."
}
